# ---------------------------------------------------------------------------
# Implementing US 25 and US 29
#  - Updated team report (Backlog "Status->Completed" column, Sprint4 sheet
#    with the new US25 / US29 rows filled in)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Backlog sheet - mark the finished stories as "Completed" (was "Complete")
# ---------------------------------------------------------------------------
$wsBacklog = $wb.Worksheets.Item("Backlog")

$completedRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,25,27,28,29,30,31,32,33)
foreach ($r in $completedRows) {
    $wsBacklog.Range("E$r").Value = "Completed"
}

# Row 19 grew a line and is now as tall as the other two-line rows
$wsBacklog.Rows.Item(19).RowHeight = 34

# Activate the Backlog sheet and update the visible selection
$wsBacklog.Activate()
$wsBacklog.Range("A2").Select()

# ---------------------------------------------------------------------------
# 2. Sprint3 sheet - no longer the tab shown when the workbook is opened
#    (handled further below when Sprint4 is activated last)
# ---------------------------------------------------------------------------
$wsSprint3 = $wb.Worksheets.Item("Sprint3")

# ---------------------------------------------------------------------------
# 3. Sprint4 sheet - fill in the newly implemented US25 and US29 rows
# ---------------------------------------------------------------------------
$wsSprint4 = $wb.Worksheets.Item("Sprint4")

# Column N is now wide enough to show the full test file name
$wsSprint4.Columns.Item(14).ColumnWidth = 17.4

# --- Row 9 : US29 "List deceased" -------------------------------------------------
$wsSprint4.Range("D9").Value = "Complete"
$wsSprint4.Range("G9").Value = 28
$wsSprint4.Range("H9").Value = 1
$wsSprint4.Range("I9").Value = 42094
$wsSprint4.Range("J9").Value = "DeseasedIndividuals.py"
$wsSprint4.Range("K9").Value = "list_deseased_individuals"
$wsSprint4.Range("O9").Value = "test_list_deseased"
$wsSprint4.Range("N9").Value = "TestDeseasedIndividauls,py"
$wsSprint4.Range("L9").Value = 20
$wsSprint4.Range("P9").Value = 20
$wsSprint4.Rows.Item(9).RowHeight = 28

# --- Row 7 : US25 "Unique first names in families" -------------------------------
$wsSprint4.Range("D7").Value = "Complete"
$wsSprint4.Range("G7").Value = 45
$wsSprint4.Range("H7").Value = 1
$wsSprint4.Range("I7").Value = 42094
$wsSprint4.Range("J7").Value = "FamilyValidation.oy"
$wsSprint4.Range("K7").Value = "check_same_name"
$wsSprint4.Range("O7").Value = "test_same_name_XXX"
$wsSprint4.Range("N7").Value = "TestFamilyValidation.py"
$wsSprint4.Range("L7").Value = 37
$wsSprint4.Range("P7").Value = 58

# Apply the date format used elsewhere in the workbook for the two new dates
$wsSprint4.Range("I7,I9").NumberFormat = "m/d/yy"

# Make Sprint4 the active (visible) sheet/tab and update its selection
$wsSprint4.Activate()
$wsSprint4.Range("P8").Select()

Write-Host "Backlog + Sprint4 updated for US25/US29"
